# Daily attendance processing - 2025-12-03 22:52:43
# Normalize the "Recorded By" (column G) lists: rotate each comma-separated
# list of recorders right by one position (the last entry moves to the front).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $val = $cell.Value2

    if ($val -ne $null -and $val -like "*,*") {
        $parts = $val -split ", "
        $count = $parts.Length
        $rotated = @($parts[$count - 1]) + $parts[0..($count - 2)]
        $newval = $rotated -join ", "
        $cell.Value = $newval
    }
}
